# Generate Report for Handoff
# - Update the "Latest Handoff Datetime" timestamps for the
#   26ad4690-800d-4a71-b1eb-9b9e2fd75956.md row group (rows 8,9,11,12,13,14)
#   on the Overview, zh-cn and de-de sheets.
# - Set the "Priority" column to "ht" for those same rows on the
#   zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows affected by this handoff report generation.
$rows = @(8, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $wsOverview.Range("G$r").Value = "2016-08-29 22:22:52"

    # de-de sheet: column H = "Latest Handoff Datetime"
    $wsDeDe.Range("H$r").Value = "2016-08-29 22:22:52"
    # de-de sheet: column E = "Priority"
    $wsDeDe.Range("E$r").Value = "ht"

    # zh-cn sheet: column H = "Latest Handoff Datetime"
    $wsZhCn.Range("H$r").Value = "2016-08-29 22:22:47"
    # zh-cn sheet: column E = "Priority"
    $wsZhCn.Range("E$r").Value = "ht"
}
